$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.821.91'
$ws.Range("E2").Value = '  +4.11%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.875.67'
$ws.Range("E3").Value = '  +3.27%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '277.21'
$ws.Range("E5").Value = '  -0.13%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.01%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5263'
$ws.Range("E7").Value = '  +3.70%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3414'
$ws.Range("E8").Value = '  -3.32%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06943'
$ws.Range("E9").Value = '  +4.17%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8034'
$ws.Range("E11").Value = '  -2.88%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07714'
$ws.Range("E12").Value = '  -1.78%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.840.49'
$ws.Range("E13").Value = '  +1.39%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.185'
$ws.Range("E14").Value = '  +2.20%  '
$ws.Range("B15").Value = 'Litecoin'
$ws.Range("C15").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '90.23'
$ws.Range("E15").Value = '  +2.92%  '
$ws.Range("E16").Value = '  +3.32%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.000'
$ws.Range("E17").Value = '  +0.01%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008047'
$ws.Range("E18").Value = '  +0.18%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.000'
$ws.Range("E19").Value = '  -0.04%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '26.865.48'
$ws.Range("E20").Value = '  +4.11%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.103.34'
$ws.Range("E21").Value = '  +2.96%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.737'
$ws.Range("E22").Value = '  -0.11%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.04'
$ws.Range("E23").Value = '  +0.33%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.174'
$ws.Range("E24").Value = '  +1.35%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.386'
$ws.Range("E25").Value = '  +8.25%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '146.74'
$ws.Range("E26").Value = '  +2.92%  '
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.34'
$ws.Range("E27").Value = '  +1.34%  '
$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.661'
$ws.Range("E28").Value = '  -0.82%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '113.42'
$ws.Range("E29").Value = '  +3.71%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.350'
$ws.Range("E30").Value = '  +0.15%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.304'
$ws.Range("E31").Value = '  +1.62%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04902'
$ws.Range("E33").Value = '  +0.62%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.173'
$ws.Range("E34").Value = '  +3.10%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7265'
$ws.Range("E35").Value = '  -0.39%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.868'
$ws.Range("E36").Value = '  -0.05%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.286'
$ws.Range("E37").Value = '  +5.01%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.339'
$ws.Range("E38").Value = '  -2.03%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01842'
$ws.Range("E39").Value = '  -0.43%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5120'
$ws.Range("E40").Value = '  -0.88%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9558'
$ws.Range("E41").Value = '  -0.88%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '116.16'
$ws.Range("E42").Value = '  +5.01%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.152'
$ws.Range("E43").Value = '  -0.88%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.110'
$ws.Range("E44").Value = '  +1.13%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.000'
$ws.Range("E45").Value = '  -0.04%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4455'
$ws.Range("E46").Value = '  -2.51%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1338'
$ws.Range("E47").Value = '  -2.02%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.313'
$ws.Range("E48").Value = '  +0.50%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '36.26'
$ws.Range("E49").Value = '  -0.97%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05944'
$ws.Range("E50").Value = '  +1.77%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.491'
$ws.Range("E51").Value = '  -0.66%  '
